$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every data cell as literal text (inline strings), even
# when the text looks like a plain number ("0.630", "71.20", ...). Excel's COM
# Range.Value setter auto-detects such strings and silently coerces them into a
# Number (dropping trailing zeros, e.g. "0.630" -> 0.63), which would corrupt the
# scraped values. Pre-format those specific cells as Text so the literal string is
# kept verbatim, matching the original data.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D15", "D16", "D20", "D21", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D46", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.896.57"
$ws.Range("E2").Value = "  -1.67%  "

$ws.Range("D3").Value = "2.231.24"
$ws.Range("E3").Value = "  -2.26%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "245.84"
$ws.Range("E5").Value = "  -2.31%  "

$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -2.59%  "

$ws.Range("D7").Value = "75.65"
$ws.Range("E7").Value = "  +2.52%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  -2.63%  "

$ws.Range("D10").Value = "40.90"
$ws.Range("E10").Value = "  +4.57%  "

$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").Value = "  -2.66%  "

$ws.Range("D12").Value = "7.08"
$ws.Range("E12").Value = "  -3.88%  "

$ws.Range("E13").Value = "  -2.35%  "

$ws.Range("D14").Value = "2.566.44"
$ws.Range("E14").Value = "  -1.99%  "

$ws.Range("D15").Value = "14.79"
$ws.Range("E15").Value = "  -1.70%  "

$ws.Range("D16").Value = "0.856"
$ws.Range("E16").Value = "  -2.56%  "

$ws.Range("D17").Value = "2.221.66"
$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("D18").Value = "41.815.98"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("D19").Value = "0.0₃0978"
$ws.Range("E19").Value = "  -2.27%  "

$ws.Range("D20").Value = "6.13"
$ws.Range("E20").Value = "  -2.71%  "

$ws.Range("D21").Value = "71.20"
$ws.Range("E21").Value = "  -1.23%  "

$ws.Range("D22").Value = "2.24"
$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("D23").Value = "230.50"
$ws.Range("E23").Value = "  -1.99%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("E25").Value = "  -4.85%  "

$ws.Range("D26").Value = "11.04"
$ws.Range("E26").Value = "  -4.16%  "

$ws.Range("D27").Value = "2.30"
$ws.Range("E27").Value = "  -5.14%  "

$ws.Range("D28").Value = "7.36"
$ws.Range("E28").Value = "  +15.24%  "

$ws.Range("D29").Value = "2.15"
$ws.Range("E29").Value = "  +0.32%  "

$ws.Range("D30").Value = "169.07"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").Value = "20.46"
$ws.Range("E31").Value = "  -3.12%  "

$ws.Range("D32").Value = "33.54"
$ws.Range("E32").Value = "  +6.20%  "

$ws.Range("D33").Value = "0.0847"
$ws.Range("E33").Value = "  +4.82%  "

$ws.Range("E34").Value = "  -5.22%  "

$ws.Range("E35").Value = "  -1.13%  "

$ws.Range("D36").Value = "4.62"
$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("D37").Value = "4.88"
$ws.Range("E37").Value = "  +2.60%  "

$ws.Range("D38").Value = "0.0299"
$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("D39").Value = "13.68"
$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Value = "5.87"
$ws.Range("E40").Value = "  -0.80%  "

$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "2.18"
$ws.Range("E41").Value = "  -7.08%  "

$ws.Range("D42").Value = "111.40"
$ws.Range("E42").Value = "  +13.67%  "

$ws.Range("E43").Value = "  -4.55%  "

$ws.Range("D44").Value = "60.19"
$ws.Range("E44").Value = "  -2.91%  "

$ws.Range("D46").Value = "0.100"
$ws.Range("E46").Value = "  -3.99%  "

$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("E48").Value = "  -4.39%  "

$ws.Range("E49").Value = "  -1.56%  "

$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "4.23"
$ws.Range("E50").Value = "  -12.83%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "2.25"
$ws.Range("E51").Value = "  -1.45%  "
